$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.611.26"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "2.518.87"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'591.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "'176.49"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.87%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "2.517.69"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +2.75%  "
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "'26.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "3.026.02"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "67.473.16"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").Value = "2.510.34"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "'7.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.20%  "
$ws.Range("D20").Value = "'11.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "'361.84"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.80%  "
$ws.Range("D22").Value = "'4.19"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'4.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'70.99"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("D27").Value = "'10.18"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.644.58"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "0.0₃0988"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "'547.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.60%  "
$ws.Range("D32").Value = "'8.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").Value = "'1.34"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").Value = "'1.86"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").Value = "'0.131"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").Value = "'18.76"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("D43").Value = "'5.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("D44").Value = "'2.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.32%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'0.561"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'146.80"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0280"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").Value = "'3.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "'0.0756"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.39%  "
